# "keck added as an admin"
#
# ProductBacklog.xlsx - row 4 ("Tutor list / Admin wants to maintain
# tutoring profiles") as well as rows 17/18 ("Contact"/"Images" stories)
# get their "Design" column (G) flipped from "no" to "yes", matching the
# look of the other already-"done" cells in the same rows. The last
# selected cell in the sheet also moved from F7 to G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - paste only the formatting (fill/font/border) of the
# copied cell, leaving the destination's current value untouched.
$xlPasteFormats = -4122

# --- G4: "no" -> "yes" -------------------------------------------------
# Re-use the formatting already applied to the other "yes" cells on row 2
# (same font/fill combination Excel uses for every "Design = yes" cell)
# then write the new value.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G4").PasteSpecial($xlPasteFormats)
$ws.Range("G4").Value = "yes"

# --- G17: "no" -> "yes" -------------------------------------------------
$ws.Range("F17").Copy() | Out-Null
$ws.Range("G17").PasteSpecial($xlPasteFormats)
$ws.Range("G17").Value = "yes"

# --- G18: "no" -> "yes" -------------------------------------------------
$ws.Range("F18").Copy() | Out-Null
$ws.Range("G18").PasteSpecial($xlPasteFormats)
$ws.Range("G18").Value = "yes"

$ws.Application.CutCopyMode = $false

# Move the active selection to G11, matching the saved view state.
$ws.Range("G11").Select()
